# Updated cryptos list on Thu May 11 06:16:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to hold a literal text value (the sheet stores prices
    # as text, e.g. "27.542.10" / "1.002" / "15.40") instead of letting
    # Excel auto-coerce a numeric-looking string into a number.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "27.542.10"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.832.56"
$ws.Range("E3").Value = "  -0.65%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.002"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
Set-TextCell "D5" "312.65"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6 - USDC (price unchanged)
$ws.Range("E6").Value = "  -0.19%  "

# Row 7 - XRP
Set-TextCell "D7" "0.4292"
$ws.Range("E7").Value = "  -0.54%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.3667"
$ws.Range("E8").Value = "  +0.64%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.07281"
$ws.Range("E9").Value = "  -0.70%  "

# Row 10 - Polygon
Set-TextCell "D10" "0.8634"
$ws.Range("E10").Value = "  -1.63%  "

# Row 11 - Solana
Set-TextCell "D11" "20.65"
$ws.Range("E11").Value = "  -0.31%  "

# Row 12 - WrappedEther
Set-TextCell "D12" "1.854.48"
$ws.Range("E12").Value = "  -1.00%  "

# Row 13 - Polkadot
Set-TextCell "D13" "5.395"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14 - Chainlink
Set-TextCell "D14" "6.521"
$ws.Range("E14").Value = "  +0.00%  "

# Row 15 - TRON (price unchanged)
$ws.Range("E15").Value = "  -0.16%  "

# Row 16 - BinanceUSD
Set-TextCell "D16" "1.002"
$ws.Range("E16").Value = "  -0.28%  "

# Row 17 - Litecoin
Set-TextCell "D17" "80.63"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18 - ShibaInu
Set-TextCell "D18" "0.000008904"
$ws.Range("E18").Value = "  -1.10%  "

# Row 19 - Dai (price unchanged)
$ws.Range("E19").Value = "  -0.28%  "

# Row 20 - Avalanche
Set-TextCell "D20" "15.40"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21 - WrappedBTC
Set-TextCell "D21" "27.768.39"
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - Uniswap
Set-TextCell "D22" "5.153"
$ws.Range("E22").Value = "  +3.48%  "

# Row 23 - Cosmos
Set-TextCell "D23" "10.83"
$ws.Range("E23").Value = "  +5.22%  "

# Row 24 - WrappedliquidstakedEther2.0
Set-TextCell "D24" "2.093.17"
$ws.Range("E24").Value = "  +1.23%  "

# Row 25 - Toncoin
Set-TextCell "D25" "1.991"
$ws.Range("E25").Value = "  +0.03%  "

# Row 26 - Monero
Set-TextCell "D26" "154.81"
$ws.Range("E26").Value = "  -0.61%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "18.86"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextCell "D28" "5.099"
$ws.Range("E28").Value = "  -2.23%  "

# Row 29 - BitcoinCash
Set-TextCell "D29" "114.29"
$ws.Range("E29").Value = "  -4.23%  "

# Row 30 - LidoDAOToken
Set-TextCell "D30" "1.823"
$ws.Range("E30").Value = "  -2.53%  "

# Row 31 - Stellar
Set-TextCell "D31" "0.08852"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32 - ImmutableX
Set-TextCell "D32" "0.7514"
$ws.Range("E32").Value = "  -0.22%  "

# Row 33 - HuobiToken
Set-TextCell "D33" "2.990"
$ws.Range("E33").Value = "  +1.10%  "

# Row 34 - Filecoin
Set-TextCell "D34" "4.542"
$ws.Range("E34").Value = "  +0.47%  "

# Row 35 - ARBITRUM
Set-TextCell "D35" "1.132"
$ws.Range("E35").Value = "  +1.28%  "

# Row 36 - Frax (price unchanged)
$ws.Range("E36").Value = "  -0.14%  "

# Row 37 - TrustWalletToken
Set-TextCell "D37" "1.090"
$ws.Range("E37").Value = "  -1.18%  "

# Row 38 - Hedera
Set-TextCell "D38" "0.05313"
$ws.Range("E38").Value = "  -2.22%  "

# Row 39 - VeChain
Set-TextCell "D39" "0.01932"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40 - MXToken
Set-TextCell "D40" "2.803"
$ws.Range("E40").Value = "  -1.08%  "

# Row 41 - TheSandbox
Set-TextCell "D41" "0.5074"
$ws.Range("E41").Value = "  -0.10%  "

# Row 42 - Algorand
Set-TextCell "D42" "0.1660"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43 - FraxShare
Set-TextCell "D43" "6.511"
$ws.Range("E43").Value = "  -1.72%  "

# Row 44 - Aptos
Set-TextCell "D44" "8.298"
$ws.Range("E44").Value = "  -0.74%  "

# Row 45 - EnergySwap
Set-TextCell "D45" "10.47"
$ws.Range("E45").Value = "  +0.77%  "

# Row 46 - was Quant, now Cronos (rows 46/47 swapped places)
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D46" "0.06487"
$ws.Range("E46").Value = "  -0.89%  "

# Row 47 - was Cronos, now Quant
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D47" "105.50"
$ws.Range("E47").Value = "  -0.07%  "

# Row 48 - Decentraland (price unchanged)
$ws.Range("E48").Value = "  +1.00%  "

# Row 49 - PaxDollar
Set-TextCell "D49" "1.000"
$ws.Range("E49").Value = "  -0.16%  "

# Row 50 - NEARProtocol
Set-TextCell "D50" "1.612"
$ws.Range("E50").Value = "  -1.64%  "

# Row 51 - Aave (price unchanged)
$ws.Range("E51").Value = "  -1.20%  "
